# "Added conditional feedback trials"
#
# Column B used to hold a fixed per-row feedback value: a "response" header
# plus "images/Incorrect.png" for the two middle (distractor) rows and
# "images/Correct.png" for the stop-signal row. Feedback is now driven
# conditionally instead, so:
#   - B1 header becomes "corrAns" (was "response")
#   - B4 (stop-signal row) becomes "space" (was "images/Correct.png")
#   - B2 and B3 (the two middle rows) are cleared entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write B4 before B1 so new shared strings are created in the same
# first-seen order ("space" then "corrAns") as the target workbook.
$ws.Range("B4").Value = "space"
$ws.Range("B1").Value = "corrAns"

# Remove the old per-row feedback-image values for the two middle rows.
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()

# Match the author's final cell selection.
[void]$ws.Range("B3").Select()
